$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Regression-test data swap (R33 preprod): replace the account number used
# by both data rows.
$ws.Range("E2").Value = 7068873718
$ws.Range("E3").Value = 7068873718

# Leave the view scrolled/selected the way the author's session ended up:
# column E visible at the left edge, cell L2 active/selected.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L2").Select()
